$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update GPA values on rows 3 and 4
$ws.Range("E3").Value = 3.54
$ws.Range("E4").Value = 3.7

# Update selection to E4 (matches final saved selection state)
$ws.Range("E4").Select()
